$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9068
$ws.Range("F3").Value = 1980
$ws.Range("F4").Value = 6634
$ws.Range("F5").Value = 179
$ws.Range("F6").Value = 2152
$ws.Range("F7").Value = 607
$ws.Range("F10").Value = 76
$ws.Range("F13").Value = 12
$ws.Range("F16").Value = 8987
$ws.Range("F25").Value = 101
$ws.Range("F27").Value = 202
$ws.Range("F28").Value = 1046
$ws.Range("F29").Value = 18
$ws.Range("F30").Value = 76
$ws.Range("F31").Value = 565
$ws.Range("F32").Value = 36
$ws.Range("F34").Value = 553
$ws.Range("F35").Value = 2380
$ws.Range("F36").Value = 884
$ws.Range("F37").Value = 555
$ws.Range("F41").Value = 308
$ws.Range("F42").Value = 185
$ws.Range("F45").Value = 30
$ws.Range("F47").Value = 22
$ws.Range("F48").Value = 4004

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 343
$ws.Range("F5").Value = 26

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 9068
$ws.Range("F5").Value = 1980
$ws.Range("F6").Value = 6634
$ws.Range("F8").Value = 607
$ws.Range("F13").Value = 76
$ws.Range("F14").Value = 26
$ws.Range("F15").Value = 12
$ws.Range("F17").Value = 8987
$ws.Range("F24").Value = 101
$ws.Range("F26").Value = 202
$ws.Range("F27").Value = 18
$ws.Range("F29").Value = 565
$ws.Range("F30").Value = 36
$ws.Range("F32").Value = 553
$ws.Range("F33").Value = 884
$ws.Range("F36").Value = 555
$ws.Range("F37").Value = 308
$ws.Range("F39").Value = 185
$ws.Range("F42").Value = 30
$ws.Range("F44").Value = 22
$ws.Range("F45").Value = 4004
